{"js": "// Update citation placeholders throughout the document body.\nconst replacements = [\n  { from: \"Ref-DJ49F2\", to: \"Ref-u214420\" },\n  { from: \"Ref-J7X2BZ\", to: \"Lee et al., 2020\" },\n  { from: \"Ref-J7X8K2\", to: \"Ref-f618755\" },\n  { from: \"Ref-AB12CD\", to: \"Ref-f913926\" }\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update citation placeholders throughout the document body.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ From = \"Ref-DJ49F2\"; To = \"Ref-u214420\" },\n    @{ From = \"Ref-J7X2BZ\"; To = \"Lee et al., 2020\" },\n    @{ From = \"Ref-J7X8K2\"; To = \"Ref-f618755\" },\n    @{ From = \"Ref-AB12CD\"; To = \"Ref-f913926\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.From\n    $find.Replacement.Text = $r.To\n    $find.Execute(\n        $r.From,    # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $r.To,      # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
